$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix species-name typos / OCR artifacts in the port-level landings table.
$ws.Range("B13").Value = "Anchovy"
$ws.Range("B14").Value = "White seabass"
$ws.Range("B23").Value = "Crab"
$ws.Range("B28").Value = "Salmon"
$ws.Range("B31").Value = "California pompano"

# Update the active selection to match the author's final view.
$ws.Range("B32").Select()
